$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old column D (Target header col) content that is no longer needed,
# and rewrite the table with the new schema: time period, Sales, Target

# New header row
$ws.Range("A1").Value = "time period"
$ws.Range("B1").Value = "Sales"
$ws.Range("C1").Value = "Target"

# Data rows
$ws.Range("A2").Value = "2019 Q2"
$ws.Range("B2").Value = "sales name A1"
$ws.Range("C2").Value = 50

$ws.Range("A3").Value = "2019 Q2"
$ws.Range("B3").Value = "sales name A2"
$ws.Range("C3").Value = 100

$ws.Range("A4").Value = "2019 Q2"
$ws.Range("B4").Value = "sales name A3"
$ws.Range("C4").Value = 300

$ws.Range("A5").Value = "2019 Q2"
$ws.Range("B5").Value = "lead name A"
$ws.Range("C5").Value = 400

# Apply the numeric style (matching the existing "Target"-style number format)
# used previously in column D, now reused for column C's numeric values.
$ws.Range("C2:C5").NumberFormat = "#,##0"
$ws.Range("C2:C5").Font.Name = "Calibri"
$ws.Range("C2:C5").Font.Size = 12

# Clear old column D entirely, it is no longer part of the table
$ws.Range("D1:D5").Clear()

# Update selection to match target state
$ws.Range("B3").Select()
